{"js": "// Move the \"FIWARE Step-by-Step\" Heading1 paragraph to the very start of\n// the document (in front of the introductory paragraph that used to lead),\n// and re-scope the two bookmarks accordingly:\n//   - \"fiware-step-by-step\" now wraps the \"FIWARE Step-by-Step\" heading text\n//   - \"how-to-use\" now starts right before the intro paragraph's text and\n//     ends right after the \"How to Use\" heading's text (i.e. it spans the\n//     intro paragraph + the \"How to Use\" heading paragraph)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst introPara = paragraphs.items[0];   // \"This is a collection of tutorials...\" (FirstParagraph)\nconst headingPara = paragraphs.items[1]; // \"FIWARE Step-by-Step\" (Heading1)\nconst howToUsePara = paragraphs.items[2]; // \"How to Use\" (Heading2)\n\nconst headingText = headingPara.text;\n\n// Drop the old bookmarks; they will be re-created in their new locations.\ncontext.document.deleteBookmark(\"fiware-step-by-step\");\ncontext.document.deleteBookmark(\"how-to-use\");\nawait context.sync();\n\n// Insert a brand-new Heading1 paragraph in front of the intro paragraph,\n// carrying the heading's text, then drop the old heading paragraph.\nconst newHeadingPara = introPara.insertParagraph(headingText, \"Before\");\nnewHeadingPara.style = \"Heading 1\";\n\nheadingPara.delete();\nawait context.sync();\n\n// Re-create \"fiware-step-by-step\" around the relocated heading text.\n// NB: use \"Content\" (not the default whole-paragraph range) so the end\n// boundary lands right after the run, inside this same paragraph, rather\n// than bleeding into the start of the following paragraph.\nconst headingRange = newHeadingPara.getRange(\"Content\");\nheadingRange.insertBookmark(\"fiware-step-by-step\");\n\n// Re-create \"how-to-use\" spanning from the very start of the intro\n// paragraph through to the very end of the \"How to Use\" heading paragraph.\nconst spanStart = introPara.getRange(\"Start\");\nconst spanEnd = howToUsePara.getRange(\"End\");\nconst howToUseRange = spanStart.expandTo(spanEnd);\nhowToUseRange.insertBookmark(\"how-to-use\");\n\nawait context.sync();\n", "ps1": "# Move the \"FIWARE Step-by-Step\" Heading1 paragraph to the very start of\n# the document (in front of the introductory paragraph that used to lead),\n# and re-scope the two bookmarks accordingly:\n#   - \"fiware-step-by-step\" now wraps the \"FIWARE Step-by-Step\" heading text\n#   - \"how-to-use\" now starts right before the intro paragraph's text and\n#     ends right after the \"How to Use\" heading's text (i.e. it spans the\n#     intro paragraph + the \"How to Use\" heading paragraph)\n\n$d = $word.ActiveDocument\n\n$headingText = \"FIWARE Step-by-Step\"\n\n# Drop the old bookmarks; they will be re-created in their new locations.\n$d.Bookmarks.Item(\"fiware-step-by-step\").Delete()\n$d.Bookmarks.Item(\"how-to-use\").Delete()\n\n# Insert a brand-new empty paragraph in front of the intro paragraph\n# (paragraph 1), give it the heading's text and Heading 1 style.\n$introRange = $d.Paragraphs(1).Range\n$introRange.InsertParagraphBefore() | Out-Null\n\n$newHeadingRange = $d.Paragraphs(1).Range\n$newHeadingRange.Text = $headingText\n$newHeadingRange.Style = \"Heading 1\"\n\n# Drop the old heading paragraph, now pushed down to index 3\n# (1 = new heading, 2 = intro, 3 = old heading, 4 = \"How to Use\").\n$d.Paragraphs(3).Range.Delete() | Out-Null\n\n# Re-create \"fiware-step-by-step\" around the relocated heading text only\n# (exclude the trailing paragraph mark).\n$headingPara = $d.Paragraphs(1).Range\n$headingBmRange = $d.Range($headingPara.Start, $headingPara.End - 1)\n$d.Bookmarks.Add(\"fiware-step-by-step\", $headingBmRange) | Out-Null\n\n# Re-create \"how-to-use\" spanning from the very start of the intro\n# paragraph through to the very end of the \"How to Use\" heading paragraph\n# (excluding its trailing paragraph mark).\n$introPara = $d.Paragraphs(2).Range\n$howToUsePara = $d.Paragraphs(3).Range\n$howToUseBmRange = $d.Range($introPara.Start, $howToUsePara.End - 1)\n$d.Bookmarks.Add(\"how-to-use\", $howToUseBmRange) | Out-Null\n"}
